$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update existing row 2 with the new "Sulawesi Barat / KELAPA / Perkebunan / PRODUKSI" record ---
# (Level, Start Year, End Year columns A, F, G are unchanged)
$ws.Range("B2").Value = "Sulawesi Barat"
$ws.Range("C2").Value = "KELAPA"
$ws.Range("D2").Value = "Perkebunan"
$ws.Range("E2").Value = "PRODUKSI"

# --- 2. Add the new merged rows (3-6) ---
$ws.Range("A3").Value = "Kabupaten"
$ws.Range("B3").Value = "Kalimantan Barat"
$ws.Range("C3").Value = "TOMAT"
$ws.Range("D3").Value = "Hortikultura"
$ws.Range("E3").Value = "PRODUKSI"
$ws.Range("F3").Value = 1970
$ws.Range("G3").Value = 2024

$ws.Range("A4").Value = "Kabupaten"
$ws.Range("B4").Value = "Kalimantan Barat"
$ws.Range("C4").Value = "TOMAT"
$ws.Range("D4").Value = "Hortikultura"
$ws.Range("E4").Value = "LUAS PANEN"
$ws.Range("F4").Value = 1970
$ws.Range("G4").Value = 2024

$ws.Range("A5").Value = "Kabupaten"
$ws.Range("B5").Value = "Kepulauan Riau"
$ws.Range("C5").Value = "KELAPA"
$ws.Range("D5").Value = "Perkebunan"
$ws.Range("E5").Value = "LUAS AREAL"
$ws.Range("F5").Value = 1970
$ws.Range("G5").Value = 2024

$ws.Range("A6").Value = "Kabupaten"
$ws.Range("B6").Value = "Sulawesi Barat"
$ws.Range("C6").Value = "KELAPA"
$ws.Range("D6").Value = "Perkebunan"
$ws.Range("E6").Value = "LUAS AREAL"
$ws.Range("F6").Value = 1970
$ws.Range("G6").Value = 2024

# Make sure new data rows carry the same (plain, non-bold, borderless) look as the
# other data rows rather than inheriting any stray formatting.
$ws.Range("A3:G6").Font.Bold = $false

# --- 3. Column widths ---
$ws.Columns.Item(1).ColumnWidth = 9.857142857142858
$ws.Columns.Item(2).ColumnWidth = 32.285714285714285
$ws.Columns.Item(3).ColumnWidth = 15.428571428571429
$ws.Columns.Item(4).ColumnWidth = 12.428571428571429
$ws.Columns.Item(5).ColumnWidth = 11.142857142857144
$ws.Columns.Item(6).ColumnWidth = 12.428571428571429
$ws.Columns.Item(7).ColumnWidth = 10.0

# --- 4. Outline level bookkeeping: the workbook keeps a cached "deepest outline level
# ever used" on the sheet even after the row that caused it is gone. Reproduce that by
# briefly grouping a row out past the used range and then removing it again. ---
$farRow = $ws.Rows.Item(100)
$farRow.OutlineLevel = 5
$farRow.EntireRow.Delete()

# --- 5. Selection as left by the editor ---
$ws.Range("B18").Select()
